$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('C2').Value = '[-, -, -, ''MEC-3B-Tec. Soldagem'']'
$ws.Range('E2').Value = '-'
$ws.Range('F2').Value = '[-, -, ''MCT-3A-Tecnologia da soldagem'', -]'
$ws.Range('C3').Value = '-'
$ws.Range('D3').Value = '-'
$ws.Range('E3').Value = '-'
$ws.Range('F3').Value = '-'
$ws.Range('B4').Value = '-'
$ws.Range('C4').Value = '-'
$ws.Range('D4').Value = '-'
$ws.Range('E4').Value = '-'
$ws.Range('C6').Value = '-'
$ws.Range('D6').Value = '-'
$ws.Range('E6').Value = '[-, ''MEC-3B-Tec. Soldagem'', -, -]'
$ws.Range('F6').Value = '[''MCT-3A-Tecnologia da soldagem'', -, -, -]'
$ws.Range('B7').Value = '-'
$ws.Range('C7').Value = '-'
$ws.Range('D7').Value = '-'
$ws.Range('E7').Value = '[-, ''MEC-3B-Tec. Soldagem'', -, -]'
$ws.Range('F7').Value = '[''MCT-3A-Tecnologia da soldagem'', -, -, -]'
$ws.Range('C8').Value = '[-, ''MEC-3B-Tec. Soldagem'', -, -]'
$ws.Range('F8').Value = '[''MCT-3A-Tecnologia da soldagem'', -, -, -]'
$ws.Range('B10').Value = '-'
$ws.Range('B11').Value = '-'
$ws.Range('F11').Value = '-'
$ws.Range('D12').Value = '-'
$ws.Range('F14').Value = '-'
$ws.Range('C15').Value = '-'
$ws.Range('D15').Value = '-'
$ws.Range('B16').Value = '-'
$ws.Range('B18').Value = '-'
$ws.Range('C18').Value = '[-, -, -, ''MEC-2NA-Soldagem'']'
$ws.Range('D18').Value = '[-, -, -, ''MEC-2NA-Soldagem'']'
$ws.Range('F18').Value = '-'
$ws.Range('B19').Value = '-'
$ws.Range('C19').Value = '[-, -, -, ''MEC-2NA-Soldagem'']'
$ws.Range('D19').Value = '[-, -, -, ''MEC-2NA-Soldagem'']'
$ws.Range('F19').Value = '-'
$ws.Range('B20').Value = 'ELM-1NA-Gestão Integrada'
$ws.Range('C20').Value = '-'
$ws.Range('D20').Value = 'MEC-2NA-Gest. Int.'
$ws.Range('E20').Value = '-'
$ws.Range('F20').Value = 'ELM-1NA-Gestão Integrada'
$ws.Range('B21').Value = 'MEC-2NB-Gestão integrada'
$ws.Range('C21').Value = '-'
$ws.Range('D21').Value = 'MEC-2NA-Gest. Int.'
$ws.Range('E21').Value = 'MEC-2NB-Gestão integrada'
